$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 1645
$ws.Range("J70").Value = 1568.875
$ws.Range("L70").Value = 4706.625
$ws.Range("N70").Value = -5246.625

$ws.Range("H73").Value = 1645
$ws.Range("J73").Value = 1568.875
$ws.Range("L73").Value = 4706.625
$ws.Range("N73").Value = -6578.625

$ws.Range("H97").Value = 900.8
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 900.8
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 2702.4
$ws.Range("N97").Value = -3694.4
$ws.Range("M97").ClearContents()

$ws.Range("H137").Value = 405128.4
$ws.Range("I137").Value = 1713.7916
$ws.Range("J137").Value = 1211957.6
$ws.Range("K137").Value = 5141.3748
$ws.Range("L137").Value = 3635872.8
$ws.Range("M137").Value = -2591.3748
$ws.Range("N137").Value = -3640972.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7605.2856
$ws.Range("J32").Value = 18947.812
$ws.Range("L32").Value = 18947.812
$ws.Range("N32").Value = -19521.812

$ws.Range("H74").Value = 53057.85
$ws.Range("I74").Value = 127303.625
$ws.Range("K74").Value = 127303.625
$ws.Range("M74").Value = -126429.625

$ws.Range("H77").Value = 53057.85
$ws.Range("I77").Value = 127303.625
$ws.Range("K77").Value = 636518.125
$ws.Range("M77").Value = -632150.125

$ws.Range("H110").Value = 826.7059
$ws.Range("I110").Value = 603.6667
$ws.Range("K110").Value = 603.6667
$ws.Range("M110").Value = 1441.3333

$ws.Range("H121").Value = 78257.60000000001
$ws.Range("J121").Value = 78257.60000000001
$ws.Range("L121").Value = 78257.60000000001
$ws.Range("N121").Value = -81751.60000000001

$ws.Range("H134").Value = 90428.5
$ws.Range("J134").Value = 90428.5
$ws.Range("L134").Value = 90428.5
$ws.Range("N134").Value = -100568.5

$ws.Range("H135").Value = 150000
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 150000
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 150000
$ws.Range("N135").Value = -160140

$ws.Range("H138").Value = 0
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("M138").ClearContents()
$ws.Range("N138").ClearContents()

$ws.Range("H139").Value = 90940.8
$ws.Range("J139").Value = 90940.8
$ws.Range("L139").Value = 90940.8
$ws.Range("N139").Value = -101220.8

$ws.Range("H141").Value = 149966.75
$ws.Range("J141").Value = 149966.75
$ws.Range("L141").Value = 149966.75
$ws.Range("N141").Value = -160326.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 128277.125
$ws.Range("J22").Value = 4498.6665
$ws.Range("L22").Value = 4498.6665
$ws.Range("N22").Value = -4844.6665

$ws.Range("H40").Value = 55000
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 55000
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 55000
$ws.Range("N40").Value = -55530

$ws.Range("H49").Value = 6250
$ws.Range("J49").Value = 6500
$ws.Range("L49").Value = 6500
$ws.Range("N49").Value = -6978

$ws.Range("H94").Value = 2330.5417
$ws.Range("I94").Value = 1655.5883
$ws.Range("J94").Value = 3969.7144
$ws.Range("K94").Value = 1655.5883
$ws.Range("L94").Value = 3969.7144
$ws.Range("M94").Value = -1204.5883
$ws.Range("N94").Value = -4871.7144

$ws.Range("H96").Value = 8952
$ws.Range("I96").Value = 8952
$ws.Range("K96").Value = 8952
$ws.Range("M96").Value = -6206

$ws.Range("H134").Value = 3479.037
$ws.Range("I134").Value = 2038.9584
$ws.Range("J134").Value = 14999.667
$ws.Range("K134").Value = 6116.8752
$ws.Range("L134").Value = 44999.001
$ws.Range("M134").Value = -3581.8752
$ws.Range("N134").Value = -50069.001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1780.1818
$ws.Range("I16").Value = 1255.5
$ws.Range("J16").Value = 1896.7778
$ws.Range("K16").Value = 1255.5
$ws.Range("L16").Value = 1896.7778
$ws.Range("M16").Value = -968.5
$ws.Range("N16").Value = -2470.7778

$ws.Range("H22").Value = 449.25
$ws.Range("I22").Value = 265.66666
$ws.Range("K22").Value = 265.66666
$ws.Range("M22").Value = 84.33334000000002

$ws.Range("H31").Value = 3536.92
$ws.Range("I31").Value = 2314
$ws.Range("K31").Value = 2314
$ws.Range("M31").Value = -2019

$ws.Range("H33").Value = 4945.25
$ws.Range("I33").Value = 1952.4
$ws.Range("J33").Value = 9933.333000000001
$ws.Range("K33").Value = 1952.4
$ws.Range("L33").Value = 9933.333000000001
$ws.Range("M33").Value = -1573.4
$ws.Range("N33").Value = -10691.333

$ws.Range("H34").Value = 3536.92
$ws.Range("I34").Value = 2314
$ws.Range("K34").Value = 2314
$ws.Range("M34").Value = -2112

$ws.Range("H69").Value = 40749.25
$ws.Range("I69").Value = 40749.25
$ws.Range("K69").Value = 40749.25
$ws.Range("M69").Value = -40000.25

$ws.Range("H70").Value = 38333.332
$ws.Range("J70").Value = 38333.332
$ws.Range("L70").Value = 38333.332
$ws.Range("N70").Value = -38963.332

$ws.Range("H72").Value = 40749.25
$ws.Range("I72").Value = 40749.25
$ws.Range("K72").Value = 122247.75
$ws.Range("M72").Value = -118503.75

$ws.Range("H73").Value = 38333.332
$ws.Range("J73").Value = 38333.332
$ws.Range("L73").Value = 38333.332
$ws.Range("N73").Value = -40517.332

$ws.Range("H99").Value = 5086671
$ws.Range("I99").Value = 4833282
$ws.Range("J99").Value = 6252260
$ws.Range("K99").Value = 4833282
$ws.Range("L99").Value = 6252260
$ws.Range("M99").Value = -4831784
$ws.Range("N99").Value = -6255256

$ws.Range("H113").Value = 1780.1818
$ws.Range("I113").Value = 1255.5
$ws.Range("J113").Value = 1896.7778
$ws.Range("K113").Value = 1255.5
$ws.Range("L113").Value = 1896.7778
$ws.Range("M113").Value = 914.5
$ws.Range("N113").Value = -6236.7778

$ws.Range("H126").Value = 5086671
$ws.Range("I126").Value = 4833282
$ws.Range("J126").Value = 6252260
$ws.Range("K126").Value = 14499846
$ws.Range("L126").Value = 18756780
$ws.Range("M126").Value = -14497376
$ws.Range("N126").Value = -18761720

$ws.Range("H141").Value = 88325
$ws.Range("I141").Value = 30000
$ws.Range("K141").Value = 30000
$ws.Range("M141").Value = -24820

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 25174492
$ws.Range("I4").Value = 33565656
$ws.Range("J4").Value = 1000
$ws.Range("K4").Value = 100696968
$ws.Range("L4").Value = 3000
$ws.Range("M4").Value = -100696856
$ws.Range("N4").Value = -3224

$ws.Range("H14").Value = 84708
$ws.Range("I14").Value = 84708
$ws.Range("K14").Value = 254124
$ws.Range("M14").Value = -253951

$ws.Range("H32").Value = 80000
$ws.Range("I32").Value = 80000
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 240000
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -239717
$ws.Range("N32").ClearContents()

$ws.Range("H59").Value = 3398.889
$ws.Range("J59").Value = 3898.2
$ws.Range("L59").Value = 11694.6
$ws.Range("N59").Value = -12774.6

$ws.Range("H125").Value = 5991.857
$ws.Range("I125").Value = 4147.6665
$ws.Range("J125").Value = 7375
$ws.Range("K125").Value = 12442.9995
$ws.Range("L125").Value = 22125
$ws.Range("M125").Value = -7522.999500000002
$ws.Range("N125").Value = -31965

$ws.Range("H126").Value = 5221.6665
$ws.Range("I126").Value = 3332.5
$ws.Range("J126").Value = 9000
$ws.Range("K126").Value = 9997.5
$ws.Range("L126").Value = 27000
$ws.Range("M126").Value = -5057.5
$ws.Range("N126").Value = -36880

$ws.Range("H132").Value = 7623.5713
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 7623.5713
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 68612.14169999999
$ws.Range("N132").Value = -73672.14169999999
$ws.Range("M132").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 3002.5

$ws.Range("H92").Value = 4999.5
$ws.Range("J92").Value = 4999.5
$ws.Range("L92").Value = 4999.5
$ws.Range("N92").Value = -8743.5

$ws.Range("H138").Value = 150000
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 150000
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 150000
$ws.Range("N138").Value = -160280

$ws.Range("H141").Value = 70498.75
$ws.Range("J141").Value = 69998.336
$ws.Range("L141").Value = 69998.336
$ws.Range("N141").Value = -80358.336

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2127.7
$ws.Range("I16").Value = 2809.8333
$ws.Range("J16").Value = 1104.5
$ws.Range("K16").Value = 2809.8333
$ws.Range("L16").Value = 1104.5
$ws.Range("M16").Value = -2639.8333
$ws.Range("N16").Value = -1444.5

$ws.Range("H22").Value = 1302.0667
$ws.Range("I22").Value = 1394.25
$ws.Range("J22").Value = 933.3333
$ws.Range("K22").Value = 1394.25
$ws.Range("L22").Value = 933.3333
$ws.Range("M22").Value = -1099.25
$ws.Range("N22").Value = -1523.3333

$ws.Range("H27").Value = 1302.0667
$ws.Range("I27").Value = 1394.25
$ws.Range("J27").Value = 933.3333
$ws.Range("K27").Value = 1394.25
$ws.Range("L27").Value = 933.3333
$ws.Range("M27").Value = -1287.25
$ws.Range("N27").Value = -1147.3333

$ws.Range("H96").Value = 0
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").ClearContents()

$ws.Range("H134").Value = 129424.5
$ws.Range("J134").Value = 129424.5
$ws.Range("L134").Value = 129424.5
$ws.Range("N134").Value = -139564.5

$ws.Range("H135").Value = 94980
$ws.Range("J135").Value = 94980
$ws.Range("L135").Value = 94980
$ws.Range("N135").Value = -105120

$ws.Range("H138").Value = 138720.5
$ws.Range("J138").Value = 138720.5
$ws.Range("L138").Value = 138720.5
$ws.Range("N138").Value = -149000.5

$ws.Range("H140").Value = 68997.8
$ws.Range("J140").Value = 69997.5
$ws.Range("L140").Value = 69997.5
$ws.Range("N140").Value = -80357.5

$ws.Range("H141").Value = 150000
$ws.Range("J141").Value = 150000
$ws.Range("L141").Value = 150000
$ws.Range("N141").Value = -160360

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("N12").ClearContents()

$ws.Range("H46").Value = 98315.836
$ws.Range("J46").Value = 98315.836
$ws.Range("L46").Value = 98315.836
$ws.Range("N46").Value = -98777.836

$ws.Range("H61").Value = 1565532.9
$ws.Range("J61").Value = 36778.5
$ws.Range("L61").Value = 36778.5
$ws.Range("N61").Value = -37362.5

$ws.Range("H81").Value = 1789
$ws.Range("I81").Value = 1789
$ws.Range("J81").Value = 1789
$ws.Range("K81").Value = 3578
$ws.Range("L81").Value = 3578
$ws.Range("M81").Value = -2517
$ws.Range("N81").Value = -5700

$ws.Range("H84").Value = 1789
$ws.Range("I84").Value = 1789
$ws.Range("J84").Value = 1789
$ws.Range("K84").Value = 17890
$ws.Range("L84").Value = 17890
$ws.Range("M84").Value = -12586
$ws.Range("N84").Value = -28498

$ws.Range("H132").Value = 1906.6666
$ws.Range("I132").Value = 1600.4828
$ws.Range("J132").Value = 2589.6924
$ws.Range("K132").Value = 4801.4484
$ws.Range("L132").Value = 7769.0772
$ws.Range("M132").Value = -2271.4484
$ws.Range("N132").Value = -12829.0772

$ws.Range("H133").Value = 50981.5
$ws.Range("J133").Value = 50981.5
$ws.Range("L133").Value = 50981.5
$ws.Range("N133").Value = -61101.5

$ws.Range("H134").Value = 98315.836
$ws.Range("J134").Value = 98315.836
$ws.Range("L134").Value = 294947.508
$ws.Range("N134").Value = -300017.508

$ws.Range("H135").Value = 89933
$ws.Range("J135").Value = 89933
$ws.Range("L135").Value = 89933
$ws.Range("N135").Value = -100073

$ws.Range("H137").Value = 150000
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 150000
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 150000
$ws.Range("N137").Value = -160200

$ws.Range("H138").Value = 150000
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 150000
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 150000
$ws.Range("N138").Value = -160280

$ws.Range("H140").Value = 54349
$ws.Range("J140").Value = 54349
$ws.Range("L140").Value = 54349
$ws.Range("N140").Value = -64709

$ws.Range("H141").Value = 150000
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 150000
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 150000
$ws.Range("N141").Value = -160360
